$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 830.75
$ws.Range("I2").Value = 308.42856
$ws.Range("K2").Value = 308.42856
$ws.Range("M2").Value = -195.42856
# Row 19
$ws.Range("H19").Value = 2400
$ws.Range("I19").Value = 2400
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2400
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -2225
$ws.Range("N19").ClearContents()
# Row 33
$ws.Range("H33").Value = 167.41667
$ws.Range("I33").Value = 169
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 169
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = -608
# Row 40
$ws.Range("H40").Value = 1620.9736
$ws.Range("I40").Value = 1521.4286
$ws.Range("J40").Value = 1899.7
$ws.Range("K40").Value = 1521.4286
$ws.Range("L40").Value = 1899.7
$ws.Range("M40").Value = -1346.4286
$ws.Range("N40").Value = -2249.7
# Row 70
$ws.Range("H70").Value = 5000
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15540
# Row 73
$ws.Range("H73").Value = 5000
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -16872
# Row 88
$ws.Range("H88").Value = 1500
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 1500
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 98
$ws.Range("H98").Value = 961.44446
$ws.Range("J98").Value = 566.6667
$ws.Range("L98").Value = 566.6667
$ws.Range("N98").Value = -3562.6667
# Row 122
$ws.Range("H122").Value = 961.44446
$ws.Range("J122").Value = 566.6667
$ws.Range("L122").Value = 1700.0001
$ws.Range("N122").Value = -6600.0001
# Row 127
$ws.Range("H127").Value = 1926.8334
$ws.Range("I127").Value = 2041.125
$ws.Range("K127").Value = 6123.375
$ws.Range("M127").Value = -1163.375

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 110
$ws.Range("H110").Value = 4198.75
$ws.Range("I110").Value = 4198.75
$ws.Range("K110").Value = 4198.75
$ws.Range("M110").Value = -2153.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
# Row 86
$ws.Range("H86").Value = 1814.8334
$ws.Range("I86").Value = 1814.8334
$ws.Range("K86").Value = 1814.8334
$ws.Range("M86").Value = -691.8334
# Row 89
$ws.Range("H89").Value = 1814.8334
$ws.Range("I89").Value = 1814.8334
$ws.Range("K89").Value = 9074.166999999999
$ws.Range("M89").Value = -3458.166999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 900002
$ws.Range("I3").Value = 900002
$ws.Range("K3").Value = 900002
$ws.Range("M3").Value = -899889
# Row 28
$ws.Range("H28").Value = 17999
$ws.Range("J28").Value = 17999
$ws.Range("L28").Value = 17999
$ws.Range("N28").Value = -18489
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 96.94444
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 152.1
$ws.Range("K2").Value = 168
$ws.Range("L2").Value = 912.5999999999999
$ws.Range("M2").Value = -55
$ws.Range("N2").Value = -1138.6
# Row 52
$ws.Range("H52").Value = 597.3333
$ws.Range("J52").Value = 597.3333
$ws.Range("L52").Value = 1791.9999
$ws.Range("N52").Value = -2323.9999
# Row 95
$ws.Range("H95").Value = 8900
$ws.Range("J95").Value = 8900
$ws.Range("L95").Value = 26700
$ws.Range("N95").Value = -30818

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 821
$ws.Range("I113").Value = 474.25
$ws.Range("J113").Value = 1283.3334
$ws.Range("K113").Value = 474.25
$ws.Range("L113").Value = 1283.3334
$ws.Range("M113").Value = 1695.75
$ws.Range("N113").Value = -5623.3334

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 764.1429000000001
$ws.Range("I22").Value = 774.8333
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 774.8333
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -479.8333
$ws.Range("N22").Value = -1290
# Row 27
$ws.Range("H27").Value = 764.1429000000001
$ws.Range("I27").Value = 774.8333
$ws.Range("J27").Value = 700
$ws.Range("K27").Value = 774.8333
$ws.Range("L27").Value = 700
$ws.Range("M27").Value = -667.8333
$ws.Range("N27").Value = -914
# Row 61
$ws.Range("H61").Value = 1128.8
$ws.Range("I61").Value = 763
$ws.Range("J61").Value = 1677.5
$ws.Range("K61").Value = 763
$ws.Range("L61").Value = 1677.5
$ws.Range("M61").Value = -561
$ws.Range("N61").Value = -2081.5
# Row 113
$ws.Range("H113").Value = 1128.8
$ws.Range("I113").Value = 763
$ws.Range("J113").Value = 1677.5
$ws.Range("K113").Value = 763
$ws.Range("L113").Value = 1677.5
$ws.Range("M113").Value = 1407
$ws.Range("N113").Value = -6017.5
# Row 122
$ws.Range("H122").Value = 5884.276
$ws.Range("J122").Value = 7536.8335
$ws.Range("L122").Value = 22610.5005
$ws.Range("N122").Value = -27510.5005
# Row 132
$ws.Range("H132").Value = 8564.125
$ws.Range("I132").Value = 10384.667
$ws.Range("J132").Value = 3102.5
$ws.Range("K132").Value = 31154.001
$ws.Range("L132").Value = 9307.5
$ws.Range("M132").Value = -28624.001
$ws.Range("N132").Value = -14367.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 11
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -2284
# Row 41
$ws.Range("H41").Value = 19648
$ws.Range("I41").Value = 19669
$ws.Range("J41").Value = 19639.6
$ws.Range("K41").Value = 19669
$ws.Range("L41").Value = 19639.6
$ws.Range("M41").Value = -19279
$ws.Range("N41").Value = -20419.6
# Row 57
$ws.Range("H57").Value = 109000
$ws.Range("I57").Value = 109000
$ws.Range("K57").Value = 109000
$ws.Range("M57").Value = -108246
# Row 74
$ws.Range("H74").Value = 21996.428
$ws.Range("J74").Value = 21903.334
$ws.Range("L74").Value = 21903.334
$ws.Range("N74").Value = -23775.334
# Row 77
$ws.Range("H77").Value = 21996.428
$ws.Range("J77").Value = 21903.334
$ws.Range("L77").Value = 65710.00199999999
$ws.Range("N77").Value = -75070.00199999999
